# Two clock-in/out punches were recorded for 2026-02-02. The previous
# "Total Duration" summary row (row 11) is replaced by the first of these
# punches, and a new row 12 is appended holding the second (still open,
# i.e. un-clocked-out) punch.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the cell values first. Prefixing date-looking text with a leading
# apostrophe keeps Excel from auto-converting it into a real date value -
# the workbook stores these as plain text, same as the other date/time
# cells already on the sheet (e.g. A9, A10, B9, C9 ...).
$ws.Range("A11").Value = "'2026-02-02"
$ws.Range("B11").Value = "17:53:08"
$ws.Range("C11").Value = "18:39:28"
$ws.Range("D11").Value = "0.77 Hours"

$ws.Range("A12").Value = "'2026-02-02"
$ws.Range("B12").Value = "21:07:41"
# C12/D12 stay empty - the punch has no clock-out time yet, so there is no
# clocked-out time or computed duration for this row.

# Copy the formatting already used by the other data rows (e.g. row 9) onto
# the new/changed cells, so the new rows look like the rest of the table.
$ws.Range("C9:D9").Copy()
$ws.Range("A11:B12").PasteSpecial(-4122)
$ws.Range("C12:D12").PasteSpecial(-4122)

$excel.CutCopyMode = 0
